$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Submission Item"
$ws.Range("B1").Value = "Definition"
